$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# The "Spannungswandler" line item (row 4) was removed from the Stückliste.
# Delete the entire row so everything below shifts up automatically,
# including styles, hyperlinks and the drawing/ink anchor.
$ws.Rows.Item(4).Delete()

# The "Spannungsverteiler" line (now still row 3) had its quantity corrected
# from 2 to 1.
$ws.Range("A3").Value = 1

# Restore a plausible selection state after the edit.
$ws.Range("B11").Select()
